# Weekly update: insert a new "Brócoli - Tercera" price record for
# Agrícola del Norte S.A. de Arica, dated 45041, right after the existing
# row for date 45012 (row 498), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 497 holds the most recent existing record (date 45012, Tercera).
# Insert a fresh blank row right after it (new row 498); this shifts the
# old rows 498:533 down to 499:534, same as Excel's Rows.Insert().
$ws.Rows.Item(498).Insert()

# Populate the new row 498 with the same record as row 497, except for a
# newer date (45041).
$ws.Cells.Item(498, 1).Value = 1
$ws.Cells.Item(498, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(498, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(498, 4).Value = 45041
$ws.Cells.Item(498, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(498, 5).Value = 15
$ws.Cells.Item(498, 6).Value = 100112023
$ws.Cells.Item(498, 7).Value = "Brócoli"
$ws.Cells.Item(498, 8).Value = "Sin especificar"
$ws.Cells.Item(498, 9).Value = "Tercera"
$ws.Cells.Item(498, 10).Value = 1300
$ws.Cells.Item(498, 11).Value = 600
$ws.Cells.Item(498, 12).Value = 700
$ws.Cells.Item(498, 13).Value = 650
$ws.Cells.Item(498, 14).Value = "$/unidad"
$ws.Cells.Item(498, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(498, 16).Value = 650
$ws.Cells.Item(498, 17).Value = 1
$ws.Cells.Item(498, 18).Value = "Hortaliza"
